# Weekly update: insert two new price records (most recent week) at the top
# of the "Femacal de La Calera - Alcachofa" data block, pushing the existing
# rows down by two positions (167-202 -> 169-204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 167, shifting everything
# below (including the old row 167) down by two rows.
$ws.Rows("167:168").Insert()

# --- New row 167: Alcachofa, Argentina(o), Primera ---
$ws.Cells.Item(167, 1).Value = 3
$ws.Cells.Item(167, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(167, 3).Value = "Coquimbo"
$ws.Cells.Item(167, 4).Value = 44476
$ws.Cells.Item(167, 5).Value = 5
$ws.Cells.Item(167, 6).Value = 100112013
$ws.Cells.Item(167, 7).Value = "Alcachofa"
$ws.Cells.Item(167, 8).Value = "Argentina(o)"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 235
$ws.Cells.Item(167, 11).Value = 8500
$ws.Cells.Item(167, 12).Value = 9000
$ws.Cells.Item(167, 13).Value = 8766
$ws.Cells.Item(167, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(167, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(167, 16).Value = 175
$ws.Cells.Item(167, 17).Value = 50
$ws.Cells.Item(167, 18).Value = "Hortaliza"

# --- New row 168: Alcachofa, Española, Extra ---
$ws.Cells.Item(168, 1).Value = 3
$ws.Cells.Item(168, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(168, 3).Value = "Coquimbo"
$ws.Cells.Item(168, 4).Value = 44476
$ws.Cells.Item(168, 5).Value = 5
$ws.Cells.Item(168, 6).Value = 100112013
$ws.Cells.Item(168, 7).Value = "Alcachofa"
$ws.Cells.Item(168, 8).Value = "Española"
$ws.Cells.Item(168, 9).Value = "Extra"
$ws.Cells.Item(168, 10).Value = 100
$ws.Cells.Item(168, 11).Value = 9000
$ws.Cells.Item(168, 12).Value = 9500
$ws.Cells.Item(168, 13).Value = 9250
$ws.Cells.Item(168, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(168, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(168, 16).Value = 308
$ws.Cells.Item(168, 17).Value = 30
$ws.Cells.Item(168, 18).Value = "Hortaliza"
